# Update leve-profit figures (currentAveragePrice / Leve price / profit columns)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets to reflect the
# latest market-board pull from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 450.30768
$ws.Range("I6").Value = 245.4
$ws.Range("J6").Value = 1133.3334
$ws.Range("K6").Value = 736.2
$ws.Range("L6").Value = 3400.0002
$ws.Range("M6").Value = -624.2
$ws.Range("N6").Value = -3624.0002
$ws.Range("H8").Value = 762.1
$ws.Range("J8").Value = 2000.3334
$ws.Range("L8").Value = 6001.0002
$ws.Range("N8").Value = -6279.0002
$ws.Range("H39").Value = 289.30768
$ws.Range("I39").Value = 276.625
$ws.Range("J39").Value = 309.6
$ws.Range("K39").Value = 829.875
$ws.Range("L39").Value = 928.8000000000001
$ws.Range("M39").Value = -533.875
$ws.Range("N39").Value = -1520.8
$ws.Range("H70").Value = 3624.4119
$ws.Range("J70").Value = 4030.3572
$ws.Range("L70").Value = 12091.0716
$ws.Range("N70").Value = -12631.0716
$ws.Range("H73").Value = 3624.4119
$ws.Range("J73").Value = 4030.3572
$ws.Range("L73").Value = 12091.0716
$ws.Range("N73").Value = -13963.0716
$ws.Range("H106").Value = 2925.3635
$ws.Range("I106").Value = 2884.6
$ws.Range("K106").Value = 2884.6
$ws.Range("M106").Value = -2253.6
$ws.Range("H113").Value = 4932.7
$ws.Range("I113").Value = 4989.7144
$ws.Range("J113").Value = 4799.6665
$ws.Range("K113").Value = 4989.7144
$ws.Range("L113").Value = 4799.6665
$ws.Range("M113").Value = -1735.7144
$ws.Range("N113").Value = -11307.6665
$ws.Range("H116").Value = 2200.3572
$ws.Range("I116").Value = 2072.182
$ws.Range("K116").Value = 2072.182
$ws.Range("M116").Value = 1369.818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2672.25
$ws.Range("I45").Value = 1908.4286
$ws.Range("J45").Value = 3741.6
$ws.Range("K45").Value = 1908.4286
$ws.Range("L45").Value = 3741.6
$ws.Range("M45").Value = -1531.4286
$ws.Range("N45").Value = -4495.6
$ws.Range("H74").Value = 3264.6667
$ws.Range("I74").Value = 2917.6
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 2917.6
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -2043.6
$ws.Range("N74").Value = -6748
$ws.Range("H77").Value = 3264.6667
$ws.Range("I77").Value = 2917.6
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 14588
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -10220
$ws.Range("N77").Value = -33736
$ws.Range("H102").Value = 5796.727
$ws.Range("I102").Value = 5553.2
$ws.Range("J102").Value = 5999.6665
$ws.Range("K102").Value = 5553.2
$ws.Range("L102").Value = 5999.6665
$ws.Range("M102").Value = -3931.2
$ws.Range("N102").Value = -9243.666499999999
$ws.Range("H122").Value = 2304.6667
$ws.Range("I122").Value = 1950
$ws.Range("J122").Value = 3014
$ws.Range("K122").Value = 5850
$ws.Range("L122").Value = 9042
$ws.Range("M122").Value = -3400
$ws.Range("N122").Value = -13942

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1672.7778
$ws.Range("I22").Value = 1672.7778
$ws.Range("K22").Value = 1672.7778
$ws.Range("M22").Value = -1499.7778
$ws.Range("H76").Value = 49999
$ws.Range("J76").Value = 49999
$ws.Range("L76").Value = 49999
$ws.Range("N76").Value = -50629
$ws.Range("H79").Value = 49999
$ws.Range("J79").Value = 49999
$ws.Range("L79").Value = 49999
$ws.Range("N79").Value = -52183
$ws.Range("H86").Value = 6806.643
$ws.Range("I86").Value = 5356.143
$ws.Range("J86").Value = 8257.143
$ws.Range("K86").Value = 5356.143
$ws.Range("L86").Value = 8257.143
$ws.Range("M86").Value = -4233.143
$ws.Range("N86").Value = -10503.143
$ws.Range("H89").Value = 6806.643
$ws.Range("I89").Value = 5356.143
$ws.Range("J89").Value = 8257.143
$ws.Range("K89").Value = 26780.715
$ws.Range("L89").Value = 41285.715
$ws.Range("M89").Value = -21164.715
$ws.Range("N89").Value = -52517.715
$ws.Range("H105").Value = 1847
$ws.Range("I105").Value = 1601.8
$ws.Range("K105").Value = 1601.8
$ws.Range("M105").Value = 145.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4599
$ws.Range("I62").Value = 5249.25
$ws.Range("K62").Value = 5249.25
$ws.Range("M62").Value = -4625.25
$ws.Range("H65").Value = 4599
$ws.Range("I65").Value = 5249.25
$ws.Range("K65").Value = 26246.25
$ws.Range("M65").Value = -23126.25
$ws.Range("H140").Value = 56373.5
$ws.Range("J140").Value = 56373.5
$ws.Range("L140").Value = 56373.5
$ws.Range("N140").Value = -66733.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 66770
$ws.Range("I4").Value = 91024.55
$ws.Range("J4").Value = 70
$ws.Range("K4").Value = 273073.65
$ws.Range("L4").Value = 210
$ws.Range("M4").Value = -272961.65
$ws.Range("N4").Value = -434
$ws.Range("H17").Value = 4783.6875
$ws.Range("J17").Value = 5095.933
$ws.Range("L17").Value = 15287.799
$ws.Range("N17").Value = -15625.799
$ws.Range("H34").Value = 672.06665
$ws.Range("J34").Value = 948.7
$ws.Range("L34").Value = 2846.1
$ws.Range("N34").Value = -3014.1
$ws.Range("H38").Value = 333.32144
$ws.Range("J38").Value = 338.33334
$ws.Range("L38").Value = 1015.00002
$ws.Range("N38").Value = -1709.00002
$ws.Range("H39").Value = 3206.1667
$ws.Range("J39").Value = 3784.4
$ws.Range("L39").Value = 11353.2
$ws.Range("N39").Value = -11941.2
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H55").Value = 5084.5835
$ws.Range("J55").Value = 5437.727
$ws.Range("L55").Value = 16313.181
$ws.Range("N55").Value = -16667.181
$ws.Range("H107").Value = 559.8
$ws.Range("I107").Value = 300
$ws.Range("J107").Value = 624.75
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 1874.25
$ws.Range("N107").Value = -5714.25
$ws.Range("M107").Value = 1020
$ws.Range("H136").Value = 8125
$ws.Range("I136").Value = 4750
$ws.Range("K136").Value = 14250
$ws.Range("M136").Value = -9150

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2938.2856
$ws.Range("J80").Value = 3326.6667
$ws.Range("L80").Value = 3326.6667
$ws.Range("N80").Value = -5322.6667
$ws.Range("H83").Value = 2938.2856
$ws.Range("J83").Value = 3326.6667
$ws.Range("L83").Value = 16633.3335
$ws.Range("N83").Value = -26617.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 990
$ws.Range("I22").Value = 990
$ws.Range("K22").Value = 990
$ws.Range("M22").Value = -695
$ws.Range("H27").Value = 990
$ws.Range("I27").Value = 990
$ws.Range("K27").Value = 990
$ws.Range("M27").Value = -883
$ws.Range("H55").Value = 1153.6364
$ws.Range("I55").Value = 1198.75
$ws.Range("K55").Value = 1198.75
$ws.Range("M55").Value = -1025.75
$ws.Range("H82").Value = 5531.25
$ws.Range("J82").Value = 5892.857
$ws.Range("L82").Value = 5892.857
$ws.Range("N82").Value = -6614.857
$ws.Range("H85").Value = 5531.25
$ws.Range("J85").Value = 5892.857
$ws.Range("L85").Value = 5892.857
$ws.Range("N85").Value = -8388.857
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11333.333
$ws.Range("J62").Value = 11333.333
$ws.Range("L62").Value = 11333.333
$ws.Range("N62").Value = -12581.333
$ws.Range("H65").Value = 11333.333
$ws.Range("J65").Value = 11333.333
$ws.Range("L65").Value = 56666.665
$ws.Range("N65").Value = -62906.665
$ws.Range("H122").Value = 3121.2
$ws.Range("I122").Value = 2301.5
$ws.Range("J122").Value = 6400
$ws.Range("K122").Value = 6904.5
$ws.Range("L122").Value = 19200
$ws.Range("M122").Value = -4454.5
$ws.Range("N122").Value = -24100
$ws.Range("H126").Value = 4838.727
$ws.Range("I126").Value = 2462.6667
$ws.Range("K126").Value = 7388.000100000001
$ws.Range("M126").Value = -4918.000100000001
$ws.Range("H132").Value = 1050.5
$ws.Range("I132").Value = 1050.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3151.5
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
